# Auto-update draw results: append the 2025-11-09 Pick 4 draw as new row 54
# (one row below the previous last data row, 53).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 54

# Leading apostrophes force these values to be stored as literal text
# (matching every other row in the sheet, where Date/Phase/Result/
# InsertedAt are all plain text), instead of being auto-coerced by Excel
# into a date serial number ("2025-11-09") or a number ("251109").
$ws.Cells.Item($row, 1).Value = "'2025-11-09"
$ws.Cells.Item($row, 2).Value = "Pick 4"
$ws.Cells.Item($row, 3).Value = "'251109"
$ws.Cells.Item($row, 4).Value = "1-6-8-6"
$ws.Cells.Item($row, 5).Value = "2025-11-09T21:36:18.636+04:00"

# Drop the implicit "Text" number format Excel applies because of the
# leading apostrophes, so the new cells keep the same default (unstyled)
# appearance used by every other row in the sheet.
$ws.Cells.Item($row, 1).ClearFormats()
$ws.Cells.Item($row, 2).ClearFormats()
$ws.Cells.Item($row, 3).ClearFormats()
$ws.Cells.Item($row, 4).ClearFormats()
$ws.Cells.Item($row, 5).ClearFormats()
